# "fixed oath refresh segments"
# Corrects stale test data in the Adobe AAM and AppNexus upload-template
# sheets, and refreshes the active sheet/selection state left over from
# the previous save.

$wb = $excel.ActiveWorkbook

# --- Adobe AAM sheet -------------------------------------------------
$wsAdobe = $wb.Worksheets.Item(2)

# Segment name for the first test segment grew into a long repeated string.
$wsAdobe.Range("B3").Value = "Test Segment 1 Test Segment 1 Test Segment 1 Test Segment 1 Test "

# Both trait folder paths were repointed at the new 2018-11-12 test trait.
$wsAdobe.Range("F3").Value = "/All Traits/TEST20181112/TEST"
$wsAdobe.Range("F4").Value = "/All Traits/TEST20181112/TEST"

# --- AppNexus sheet ----------------------------------------------------
$wsAppNexus = $wb.Worksheets.Item(3)

# Report email no longer CCs ykoh@eyeota.com.
$wsAppNexus.Range("M3").Value = "asoh@eyeota.com"

# --- Restore / update the view state ------------------------------------
# AppNexus selection moved to C3.
$wsAppNexus.Activate()
$wsAppNexus.Range("C3").Select()

# Adobe AAM becomes the active sheet with F5 selected.
$wsAdobe.Activate()
$wsAdobe.Range("F5").Select()
